$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header in C1 from "Bank" to "Bank account"
$ws.Range("C1").Value = "Bank account"

# Update the last active selection to match the saved view state
$ws.Range("I15").Select()
